# Fixed UpdatedVinRenewal Tests to use valid VIN data and reset to original
# values afterwards.
#
# The workbook holds a small VIN lookup table on Sheet1. Rows 2-5 all share
# the same (previously invalid) VIN in column A and the same MAKE_TEXT value
# in column E for row 2. Update them to the new valid VIN / MAKE_TEXT values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VIN column (A2:A5) -> new valid VIN
$ws.Range("A2:A5").Value = "5TFUY5F1&D"

# MAKE_TEXT (E2) -> updated make text
$ws.Range("E2").Value = "TOYOTA_UPDATED"

# Restore the last active selection to F11
$ws.Range("F11").Select()
